# Fix the "Handcrafted Lightbow" card name typo and clear an erroneous
# affiliation value that had been left in F66.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reprint_list")

# Correct the card name in E89 ("Handcrafted Lightbow" -> "Handcrafted Light Bow")
$ws.Range("E89").Value = "Handcrafted Light Bow"

# Clear the stray value in F66 (was "Aw")
$ws.Range("F66").Value = ""

# Reflect the active cell left by the edit session
$ws.Activate()
$ws.Range("E90").Select()
